$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 59
$ws.Cells.Item(2, 3).Value = "face/face012.jpg"
$ws.Cells.Item(2, 4).Value = "mieten"
$ws.Cells.Item(2, 5).Value = "face"

$ws.Cells.Item(3, 2).Value = 92
$ws.Cells.Item(3, 3).Value = "face/face013.jpg"
$ws.Cells.Item(3, 4).Value = "scheitern"
$ws.Cells.Item(3, 5).Value = "face"

$ws.Cells.Item(4, 2).Value = 56
$ws.Cells.Item(4, 3).Value = "house/house015.jpg"
$ws.Cells.Item(4, 4).Value = "bitten"
$ws.Cells.Item(4, 5).Value = "house"

$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(5, 3).Value = "house/house005.jpg"
$ws.Cells.Item(5, 4).Value = "laufen"
$ws.Cells.Item(5, 5).Value = "house"

$ws.Cells.Item(6, 2).Value = 110
$ws.Cells.Item(6, 3).Value = "face/face001.jpg"
$ws.Cells.Item(6, 4).Value = "runden"
$ws.Cells.Item(6, 5).Value = "face"

$ws.Cells.Item(7, 2).Value = 11
$ws.Cells.Item(7, 3).Value = "house/house020.jpg"
$ws.Cells.Item(7, 4).Value = "schätzen"
$ws.Cells.Item(7, 5).Value = "house"

$ws.Cells.Item(8, 2).Value = 94
$ws.Cells.Item(8, 3).Value = "face/face029.jpg"
$ws.Cells.Item(8, 4).Value = "drehen"
$ws.Cells.Item(8, 5).Value = "face"

$ws.Cells.Item(9, 2).Value = 70
$ws.Cells.Item(9, 3).Value = "face/face011.jpg"
$ws.Cells.Item(9, 4).Value = "töten"
$ws.Cells.Item(9, 5).Value = "face"

$ws.Cells.Item(10, 2).Value = 91
$ws.Cells.Item(10, 3).Value = "face/face015.jpg"
$ws.Cells.Item(10, 4).Value = "nehmen"
$ws.Cells.Item(10, 5).Value = "face"

$ws.Cells.Item(11, 2).Value = 111
$ws.Cells.Item(11, 3).Value = "house/house007.jpg"
$ws.Cells.Item(11, 4).Value = "schenken"
$ws.Cells.Item(11, 5).Value = "house"

$ws.Cells.Item(12, 2).Value = 38
$ws.Cells.Item(12, 3).Value = "face/face022.jpg"
$ws.Cells.Item(12, 4).Value = "schmecken"
$ws.Cells.Item(12, 5).Value = "face"

$ws.Cells.Item(13, 2).Value = 67
$ws.Cells.Item(13, 3).Value = "house/house013.jpg"
$ws.Cells.Item(13, 4).Value = "husten"
$ws.Cells.Item(13, 5).Value = "house"

$ws.Cells.Item(14, 2).Value = 51
$ws.Cells.Item(14, 3).Value = "house/house027.jpg"
$ws.Cells.Item(14, 4).Value = "bleiben"
$ws.Cells.Item(14, 5).Value = "house"

$ws.Cells.Item(15, 2).Value = 17
$ws.Cells.Item(15, 3).Value = "house/house010.jpg"
$ws.Cells.Item(15, 4).Value = "opfern"
$ws.Cells.Item(15, 5).Value = "house"

$ws.Cells.Item(16, 2).Value = 100
$ws.Cells.Item(16, 3).Value = "face/face030.jpg"
$ws.Cells.Item(16, 4).Value = "kaufen"
$ws.Cells.Item(16, 5).Value = "face"

$ws.Cells.Item(17, 2).Value = 84
$ws.Cells.Item(17, 3).Value = "face/face016.jpg"
$ws.Cells.Item(17, 4).Value = "währen"
$ws.Cells.Item(17, 5).Value = "face"

$ws.Cells.Item(18, 2).Value = 7
$ws.Cells.Item(18, 3).Value = "house/house024.jpg"
$ws.Cells.Item(18, 4).Value = "wiegen"
$ws.Cells.Item(18, 5).Value = "house"

$ws.Cells.Item(19, 2).Value = 32
$ws.Cells.Item(19, 3).Value = "face/face009.jpg"
$ws.Cells.Item(19, 4).Value = "füttern"
$ws.Cells.Item(19, 5).Value = "face"

$ws.Cells.Item(20, 2).Value = 118
$ws.Cells.Item(20, 3).Value = "house/house001.jpg"
$ws.Cells.Item(20, 4).Value = "klappen"
$ws.Cells.Item(20, 5).Value = "house"

$ws.Cells.Item(21, 2).Value = 50
$ws.Cells.Item(21, 3).Value = "face/face005.jpg"
$ws.Cells.Item(21, 4).Value = "drohen"
$ws.Cells.Item(21, 5).Value = "face"

$ws.Cells.Item(22, 2).Value = 43
$ws.Cells.Item(22, 3).Value = "house/house021.jpg"
$ws.Cells.Item(22, 4).Value = "posten"
$ws.Cells.Item(22, 5).Value = "house"

$ws.Cells.Item(23, 2).Value = 105
$ws.Cells.Item(23, 3).Value = "face/face025.jpg"
$ws.Cells.Item(23, 4).Value = "liefern"
$ws.Cells.Item(23, 5).Value = "face"

$ws.Cells.Item(24, 2).Value = 22
$ws.Cells.Item(24, 3).Value = "house/house016.jpg"
$ws.Cells.Item(24, 4).Value = "raten"
$ws.Cells.Item(24, 5).Value = "house"

$ws.Cells.Item(25, 2).Value = 54
$ws.Cells.Item(25, 3).Value = "house/house000.jpg"
$ws.Cells.Item(25, 4).Value = "loben"
$ws.Cells.Item(25, 5).Value = "house"

$ws.Cells.Item(26, 2).Value = 44
$ws.Cells.Item(26, 3).Value = "face/face027.jpg"
$ws.Cells.Item(26, 4).Value = "sieben"
$ws.Cells.Item(26, 5).Value = "face"

$ws.Cells.Item(27, 2).Value = 2
$ws.Cells.Item(27, 3).Value = "house/house009.jpg"
$ws.Cells.Item(27, 4).Value = "formen"
$ws.Cells.Item(27, 5).Value = "house"

$ws.Cells.Item(28, 2).Value = 19
$ws.Cells.Item(28, 3).Value = "face/face019.jpg"
$ws.Cells.Item(28, 4).Value = "haken"
$ws.Cells.Item(28, 5).Value = "face"

$ws.Cells.Item(29, 2).Value = 13
$ws.Cells.Item(29, 3).Value = "face/face020.jpg"
$ws.Cells.Item(29, 4).Value = "spielen"
$ws.Cells.Item(29, 5).Value = "face"

$ws.Cells.Item(30, 2).Value = 115
$ws.Cells.Item(30, 3).Value = "house/house014.jpg"
$ws.Cells.Item(30, 4).Value = "jubeln"
$ws.Cells.Item(30, 5).Value = "house"

$ws.Cells.Item(31, 2).Value = 24
$ws.Cells.Item(31, 3).Value = "house/house025.jpg"
$ws.Cells.Item(31, 4).Value = "backen"
$ws.Cells.Item(31, 5).Value = "house"

$ws.Cells.Item(32, 2).Value = 18
$ws.Cells.Item(32, 3).Value = "face/face010.jpg"
$ws.Cells.Item(32, 4).Value = "wenden"
$ws.Cells.Item(32, 5).Value = "face"

$ws.Cells.Item(33, 2).Value = 93
$ws.Cells.Item(33, 3).Value = "house/house003.jpg"
$ws.Cells.Item(33, 4).Value = "ehren"
$ws.Cells.Item(33, 5).Value = "house"
